# gecko issue solved - 3rd Commit
#
# The underlying content change is a one-word label fix on Sheet1:
# the header in A1 ("BrowserName") had a typo and is corrected to
# "BrowseName". Everything else in the diff (shared-string table
# reshuffling, workbook window size, selection rectangle, column
# widths, page orientation) is cosmetic view-state that Excel rewrites
# whenever the file is re-saved; we reproduce the cosmetic bits too,
# on a best-effort basis, via the corresponding object-model calls.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- the actual data edit -------------------------------------------------
# Fix the mislabeled header "BrowserName" -> "BrowseName" (A1 on Sheet1).
$ws1.Range("A1").Value = "BrowseName"

# --- cosmetic / view-state edits (best effort) ----------------------------

# Selection moved from B5 to B7.
[void]$ws1.Range("B7").Select()

# Column A narrows slightly (content got one character shorter) and
# column C widens - reproduce the resulting column widths.
$ws1.Columns.Item(1).ColumnWidth = 12.0221354166667
$ws1.Columns.Item(3).ColumnWidth = 14.1666666666667

# Page orientation explicitly set to portrait (adds <pageSetup .../>).
$ws1.PageSetup.Orientation = 1

# Resize the workbook window (best effort - matches the bookViews change).
$win = $wb.Windows.Item(1)
$win.Width  = 11625
$win.Height = 4980
